$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells, matching the existing header row's style (same as AC1)
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record (Wins/Losses/Ties) for every data row
for ($r = 2; $r -le 50; $r++) {
    $ws.Cells.Item($r, 30).Value = 68   # AD: Wins
    $ws.Cells.Item($r, 31).Value = 94   # AE: Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF: Ties
}
